$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.872.91"
$ws.Range("E2").Value = "  -2.07%  "
$ws.Range("D3").Value = "3.231.00"
$ws.Range("E3").Value = "  -1.50%  "
$ws.Range("E4").Value = "  -0.01%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "578.16"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.09%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "173.57"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -3.19%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.631"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +0.83%  "
$ws.Range("D9").Value = "3.230.20"
$ws.Range("E9").Value = "  -1.39%  "
$ws.Range("E10").Value = "  -2.78%  "
$ws.Range("E11").Value = "  +1.15%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.389"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -2.90%  "
$ws.Range("D13").Value = "3.790.85"
$ws.Range("E13").Value = "  -1.50%  "
$ws.Range("E14").Value = "  -3.13%  "
$ws.Range("D15").Value = "64.944.72"
$ws.Range("E15").Value = "  -1.96%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "25.65"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -2.71%  "
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.0000159"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -2.09%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.212.90"
$ws.Range("E18").Value = "  -1.74%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "416.43"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -4.41%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "12.83"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -1.84%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "7.19"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -2.46%  "
$ws.Range("E23").Value = "  +0.09%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "70.37"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -1.93%  "
$ws.Range("E25").Value = "  -0.38%  "
$ws.Range("E26").Value = "  +3.65%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "0.0000110"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -1.94%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "9.14"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +3.21%  "
$ws.Range("E30").Value = "  +0.15%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "1.88"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -2.87%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "21.79"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -1.88%  "
$ws.Range("E33").Value = "  +0.05%  "
$ws.Range("E34").Value = "  -2.81%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "6.42"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -2.21%  "
$ws.Range("E36").Value = "  -1.85%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "156.92"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -0.45%  "
$ws.Range("E38").Value = "  -1.71%  "
$ws.Range("D39").Value = "2.827.45"
$ws.Range("E39").Value = "  +2.29%  "
$ws.Range("E40").Value = "  -2.20%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "25.46"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -4.39%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "4.23"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -1.44%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.728"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -5.72%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "39.53"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -1.81%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "5.77"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -4.23%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.0629"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -4.02%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "305.94"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -4.79%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "2.19"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -4.63%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "22.10"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -4.60%  "
$ws.Range("E50").Value = "  -0.92%  "
$ws.Range("E51").Value = "  -0.47%  "
